$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.902.13'
$ws.Range('E2').Value = '  -0.54%  '
$ws.Range('D3').Value = '2.351.42'
$ws.Range('E3').Value = '  -1.19%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('B5').Value = 'XRP'
$ws.Range('C5').Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range('D5').Value = '''0.671'
$ws.Range('E5').Value = '  -3.48%  '
$ws.Range('B6').Value = 'BNB'
$ws.Range('C6').Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range('D6').Value = '''239.83'
$ws.Range('E6').Value = '  -1.54%  '
$ws.Range('D7').Value = '''72.29'
$ws.Range('E7').Value = '  -6.87%  '
$ws.Range('D9').Value = '''0.593'
$ws.Range('E9').Value = '  -2.88%  '
$ws.Range('D10').Value = '''0.100'
$ws.Range('E10').Value = '  -4.58%  '
$ws.Range('D11').Value = '''58.36'
$ws.Range('E11').Value = '  +0.83%  '
$ws.Range('D12').Value = '''32.72'
$ws.Range('E12').Value = '  +0.12%  '
$ws.Range('E13').Value = '  -0.44%  '
$ws.Range('D14').Value = '''7.23'
$ws.Range('E14').Value = '  -4.32%  '
$ws.Range('D15').Value = '2.697.24'
$ws.Range('E15').Value = '  -1.22%  '
$ws.Range('D16').Value = '''16.32'
$ws.Range('E16').Value = '  -5.72%  '
$ws.Range('D17').Value = '''0.904'
$ws.Range('E17').Value = '  -2.90%  '
$ws.Range('D18').Value = '2.348.83'
$ws.Range('E18').Value = '  -1.17%  '
$ws.Range('D19').Value = '43.819.88'
$ws.Range('E19').Value = '  -1.27%  '
$ws.Range('D20').Value = '''0.0000103'
$ws.Range('E20').Value = '  -1.88%  '
$ws.Range('D21').Value = '''6.71'
$ws.Range('E21').Value = '  -0.61%  '
$ws.Range('D22').Value = '''78.41'
$ws.Range('E22').Value = '  -0.70%  '
$ws.Range('D23').Value = '''254.05'
$ws.Range('E23').Value = '  -1.74%  '
$ws.Range('D24').Value = '''1.93'
$ws.Range('E24').Value = '  +6.74%  '
$ws.Range('E25').Value = '  +0.04%  '
$ws.Range('D26').Value = '''3.73'
$ws.Range('E26').Value = '  +0.46%  '
$ws.Range('E27').Value = '  -2.66%  '
$ws.Range('D28').Value = '''10.47'
$ws.Range('E28').Value = '  -4.76%  '
$ws.Range('E29').Value = '  -1.29%  '
$ws.Range('D30').Value = '''177.43'
$ws.Range('E30').Value = '  +0.96%  '
$ws.Range('D31').Value = '''22.40'
$ws.Range('E31').Value = '  -3.48%  '
$ws.Range('D32').Value = '''0.128'
$ws.Range('E32').Value = '  -2.61%  '
$ws.Range('D34').Value = '''0.0746'
$ws.Range('E34').Value = '  -2.38%  '
$ws.Range('D35').Value = '''5.13'
$ws.Range('E35').Value = '  -5.04%  '
$ws.Range('D36').Value = '''5.41'
$ws.Range('E36').Value = '  +0.73%  '
$ws.Range('D37').Value = '''3.75'
$ws.Range('E37').Value = '  -4.30%  '
$ws.Range('D38').Value = '''6.42'
$ws.Range('E38').Value = '  -2.89%  '
$ws.Range('D39').Value = '''2.38'
$ws.Range('E39').Value = '  -5.08%  '
$ws.Range('D40').Value = '''0.0275'
$ws.Range('E40').Value = '  -1.62%  '
$ws.Range('D41').Value = '''66.13'
$ws.Range('E41').Value = '  +20.15%  '
$ws.Range('D42').Value = '''5.21'
$ws.Range('E42').Value = '  +15.74%  '
$ws.Range('B43').Value = 'Cronos'
$ws.Range('C43').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D43').Value = '''0.109'
$ws.Range('E43').Value = '  +7.21%  '
$ws.Range('B44').Value = 'FraxShare'
$ws.Range('C44').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D44').Value = '''9.23'
$ws.Range('E44').Value = '  +0.28%  '
$ws.Range('D45').Value = '''18.78'
$ws.Range('E45').Value = '  -2.24%  '
$ws.Range('D46').Value = '''0.199'
$ws.Range('E46').Value = '  +0.78%  '
$ws.Range('E47').Value = '  -0.10%  '
$ws.Range('B48').Value = 'NEARProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D48').Value = '''2.47'
$ws.Range('E48').Value = '  -3.94%  '
$ws.Range('B49').Value = 'TrustWalletToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D49').Value = '''1.24'
$ws.Range('E49').Value = '  -3.49%  '
$ws.Range('D50').Value = '''99.29'
$ws.Range('E50').Value = '  -3.83%  '
$ws.Range('E51').Value = '  -6.20%  '
